# Update countries & provincias Spain
# Applies the 9-Sep-2020 10:37 data refresh to the "Pais" sheet:
#   - Updated case statistics for several countries (rows whose totals changed)
#   - Eslovaquia's total cases overtook Nicaragua & Republica de Africa Central,
#     so it now sorts above them (rows 116-118 shuffle)
#   - Estonia's total cases overtook Jordania, so it now sorts above it
#     (rows 138-139 shuffle)
#   - Refreshed "Datos actualizados ..." timestamp cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Timestamp banner (row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 9 de Septiembre de 2020 a las 10:37"

# --- Straightforward statistic refreshes (country order unchanged) ---
# India (row 5)
Set-Row 5 4374314 6878 3398844 901517 0 30 73953

# Rusia (row 7)
Set-Row 7 1041007 5218 856458 166414 0 142 18135

# Israel (row 28)
Set-Row 28 139013 1448 107599 30366 0 8 1048

# Polonia (row 49)
Set-Row 49 71947 421 55910 13890 0 11 2147

# Singapur (row 53) - only totals/new-cases/recovered changed
$ws.Cells.Item(53, 2).Value = 57166
$ws.Cells.Item(53, 3).Value = 75
$ws.Cells.Item(53, 5).Value = 678

# Croacia (row 90)
Set-Row 90 12626 341 9833 2587 0 3 206

# --- Reorder block: Nicaragua / Republica de Africa Central / Eslovaquia ---
# Eslovaquia's updated total (4888) now outranks the other two (unchanged),
# so it moves to the top of the trio while the other two shift down one row.
$ws.Cells.Item(116, 1).Value = "Eslovaquia"
Set-Row 116 4888 161 2947 1904 0 0 37
$ws.Cells.Item(117, 1).Value = "Nicaragua"
Set-Row 117 4818 0 2913 1761 0 0 144
$ws.Cells.Item(118, 1).Value = "Republica de Africa Central"
Set-Row 118 4735 0 1825 2848 0 0 62

# --- Reorder block: Jordania / Estonia ---
# Estonia's updated total (2585) now outranks Jordania (unchanged),
# so it moves above it.
$ws.Cells.Item(138, 1).Value = "Estonia"
Set-Row 138 2585 22 2213 308 0 0 64
$ws.Cells.Item(139, 1).Value = "Jordania"
Set-Row 139 2581 0 1885 677 0 0 19
